$wb = $excel.ActiveWorkbook

# =========================================================
# Sheet "header" (sheet1)
# =========================================================
$ws1 = $wb.Worksheets.Item(1)

# Update the curie_map "estuarine" line (previously a TBA placeholder)
$ws1.Range("A3").Value = "   estuarine: https://w3id.org/env/neap/estuarine/"

# Insert a new line for the sssom curie, just before the status curie line
$ws1.Rows.Item(8).Insert()
$ws1.Range("A8").Value = "   sssom: https://w3id.org/sssom/"

# =========================================================
# Sheet "SSSOM" (sheet2)
# =========================================================
$ws2 = $wb.Worksheets.Item(2)

# Insert a new top header row using the full RDF/SSSOM predicate names
$ws2.Rows.Item(1).Insert()

$ws2.Range("A1").Value = "rdf:subject"
$ws2.Range("B1").Value = "sssom:subject_label"
$ws2.Range("C1").Value = "rdf:predicate"
$ws2.Range("D1").Value = "rdf:object"
$ws2.Range("E1").Value = "sssom:object_label"
$ws2.Range("F1").Value = "sssom:mapping_justification"
$ws2.Range("G1").Value = "dcterms:creator"
$ws2.Range("H1").Value = "sssom:creator_label"
$ws2.Range("I1").Value = "dcterms:created"
$ws2.Range("J1").Value = "sssom:confidence"
$ws2.Range("K1").Value = "crosswalk:status"
$ws2.Range("M1").Value = "rdfs:comment"
$ws2.Range("L1").Value = "sssom:reviewer_id"
$ws2.Range("N1").Value = "rdfs:label"

# Row-level + cell-level formatting for the new header row
$ws2.Rows.Item(1).RowHeight = 12.75
$ws2.Rows.Item(1).Font.Bold = $true
$ws2.Rows.Item(1).Font.Size = 10
$ws2.Rows.Item(1).Font.Name = "Calibri"

$hdrRange = $ws2.Range("A1:N1")
$hdrRange.Font.Bold = $true
$hdrRange.Font.Size = 9
$hdrRange.Font.Name = "Calibri"
$hdrRange.HorizontalAlignment = -4131
$hdrRange.VerticalAlignment = -4108

# Add the "mapping to IUCN GET" label formulas for the two data rows
$ws2.Range("N3").Formula = '=CONCAT(A3, " - mapping to IUCN GET")'
$ws2.Range("N4").Formula = '=CONCAT(A4, " - mapping to IUCN GET")'

# Extend the repeated blank-styled filler rows by one row (31 -> 32)
$ws2.Range("I31").NumberFormat = $ws2.Range("I30").NumberFormat
$ws2.Range("G32").Value = $ws2.Range("G31").Value

# =========================================================
# View-state (selection) restore
# =========================================================
$ws1.Range("A8:XFD8").Select()
$ws2.Activate()
$ws2.Range("M1").Select()
